# Weekly data refresh: insert one new observation row at row 31 (Albahaca /
# Femacal de La Calera), pushing the existing rows 31-130 down to 32-131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 31, shifting rows 31-130 down
# to 32-131 (this is what grows the sheet's used range from R130 to R131).
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record's data.
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "Femacal de La Calera"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44592
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 100112052
$ws.Range("G31").Value = "Albahaca"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 105
$ws.Range("K31").Value = 4000
$ws.Range("L31").Value = 4500
$ws.Range("M31").Value = 4238
$ws.Range("N31").Value = "$/docena de matas"
$ws.Range("O31").Value = "Provincia de Quillota"
$ws.Range("P31").Value = 706
$ws.Range("Q31").Value = 6
$ws.Range("R31").Value = "Hortaliza"
